# Update cell C10 on the active sheet ("Rules") from 18 to 1.
# This is the "From" value for rule R30 in the lookup table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
